$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2994946666666667
$ws.Range("H2").Value = 0.8984840000000001
$ws.Range("I2").Value = 0.4989451716962827
$ws.Range("J2").Value = 0.4989451716962828
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.951641
$ws.Range("N2").Value = 2.854923
$ws.Range("O2").Value = 0.1812272686155736
$ws.Range("P2").Value = 0.1812272686155736
$ws.Range("Q2").Value = 0.2850114040813334
$ws.Range("R2").Value = 2.565102636732
$ws.Range("S2").Value = 0.09042247065544572
$ws.Range("T2").Value = 0.09042247065544572

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2994946666666667
$ws.Range("H3").Value = 0.8984840000000001
$ws.Range("I3").Value = 0.4989451716962827
$ws.Range("J3").Value = 0.4989451716962828
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.566336
$ws.Range("N3").Value = 7.699008000000001
$ws.Range("O3").Value = 0.4887242811415405
$ws.Range("P3").Value = 0.4887242811415405
$ws.Range("Q3").Value = 0.7686039448746668
$ws.Range("R3").Value = 6.917435503872001
$ws.Range("S3").Value = 0.2438466203663083
$ws.Range("T3").Value = 0.2438466203663083

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2994946666666667
$ws.Range("H4").Value = 0.8984840000000001
$ws.Range("I4").Value = 0.4989451716962827
$ws.Range("J4").Value = 0.4989451716962828
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.733114666666667
$ws.Range("N4").Value = 5.199344
$ws.Range("O4").Value = 0.330048450242886
$ws.Range("P4").Value = 0.330048450242886
$ws.Range("Q4").Value = 0.5190585993884445
$ws.Range("R4").Value = 4.671527394496001
$ws.Range("S4").Value = 0.1646760806745288
$ws.Range("T4").Value = 0.1646760806745288

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.300761
$ws.Range("H5").Value = 0.9022829999999999
$ws.Range("I5").Value = 0.5010548283037172
$ws.Range("J5").Value = 0.5010548283037172
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.951641
$ws.Range("N5").Value = 2.854923
$ws.Range("O5").Value = 0.1812272686155736
$ws.Range("P5").Value = 0.1812272686155736
$ws.Range("Q5").Value = 0.286216498801
$ws.Range("R5").Value = 2.575948489209
$ws.Range("S5").Value = 0.09080479796012786
$ws.Range("T5").Value = 0.09080479796012786

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.300761
$ws.Range("H6").Value = 0.9022829999999999
$ws.Range("I6").Value = 0.5010548283037172
$ws.Range("J6").Value = 0.5010548283037172
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.566336
$ws.Range("N6").Value = 7.699008000000001
$ws.Range("O6").Value = 0.4887242811415405
$ws.Range("P6").Value = 0.4887242811415405
$ws.Range("Q6").Value = 0.7718537816960001
$ws.Range("R6").Value = 6.946684035264
$ws.Range("S6").Value = 0.2448776607752322
$ws.Range("T6").Value = 0.2448776607752322

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.2994946666666667
$ws.Range("H7").Value = 0.9022829999999999
$ws.Range("I7").Value = 0.5010548283037172
$ws.Range("J7").Value = 0.5010548283037172
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.733114666666667
$ws.Range("N7").Value = 5.199344
$ws.Range("O7").Value = 0.330048450242886
$ws.Range("P7").Value = 0.330048450242886
$ws.Range("Q7").Value = 0.5212533002613333
$ws.Range("R7").Value = 4.691279702351999
$ws.Range("S7").Value = 0.1653723695683572
$ws.Range("T7").Value = 0.1653723695683572
